$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 338, pushing the old rows
# 338-342 down to 341-345 (their contents/styles shift automatically).
$ws.Rows("338:340").Insert()

# --- New row 338: Mandarina / Clementina / Primera ---
$ws.Range("A338").Value = 10
$ws.Range("B338").Value = "Vega Modelo de Temuco"
$ws.Range("C338").Value = "La Araucanía"
$ws.Range("D338").Value = 44448
$ws.Range("E338").Value = 9
$ws.Range("F338").Value = "Fruta"
$ws.Range("G338").Value = 100102
$ws.Range("H338").Value = "Cítricos"
$ws.Range("I338").Value = 100102004
$ws.Range("J338").Value = "Mandarina"
$ws.Range("K338").Value = "Clementina"
$ws.Range("L338").Value = "Primera"
$ws.Range("M338").Value = 20
$ws.Range("N338").Value = 157000
$ws.Range("O338").Value = 157000
$ws.Range("P338").Value = 157000
$ws.Range("Q338").Value = "`$/bins (450 kilos)"
$ws.Range("R338").Value = "Región de O'Higgins"
$ws.Range("S338").Value = 349
$ws.Range("T338").Value = 450

# --- New row 339: Mandarina / Clementina / Tercera ---
$ws.Range("A339").Value = 10
$ws.Range("B339").Value = "Vega Modelo de Temuco"
$ws.Range("C339").Value = "La Araucanía"
$ws.Range("D339").Value = 44448
$ws.Range("E339").Value = 9
$ws.Range("F339").Value = "Fruta"
$ws.Range("G339").Value = 100102
$ws.Range("H339").Value = "Cítricos"
$ws.Range("I339").Value = 100102004
$ws.Range("J339").Value = "Mandarina"
$ws.Range("K339").Value = "Clementina"
$ws.Range("L339").Value = "Tercera"
$ws.Range("M339").Value = 7
$ws.Range("N339").Value = 100000
$ws.Range("O339").Value = 100000
$ws.Range("P339").Value = 100000
$ws.Range("Q339").Value = "`$/bins (450 kilos)"
$ws.Range("R339").Value = "Región de O'Higgins"
$ws.Range("S339").Value = 222
$ws.Range("T339").Value = 450

# --- New row 340: Mandarina / Murcott / Primera ---
$ws.Range("A340").Value = 10
$ws.Range("B340").Value = "Vega Modelo de Temuco"
$ws.Range("C340").Value = "La Araucanía"
$ws.Range("D340").Value = 44448
$ws.Range("E340").Value = 9
$ws.Range("F340").Value = "Fruta"
$ws.Range("G340").Value = 100102
$ws.Range("H340").Value = "Cítricos"
$ws.Range("I340").Value = 100102004
$ws.Range("J340").Value = "Mandarina"
$ws.Range("K340").Value = "Murcott"
$ws.Range("L340").Value = "Primera"
$ws.Range("M340").Value = 300
$ws.Range("N340").Value = 7000
$ws.Range("O340").Value = 7000
$ws.Range("P340").Value = 7000
$ws.Range("Q340").Value = "`$/bandeja 18 kilos"
$ws.Range("R340").Value = "Región de O'Higgins"
$ws.Range("S340").Value = 389
$ws.Range("T340").Value = 18
